$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 17239.166
$ws.Range("I12").Value = 17239.166
$ws.Range("K12").Value = 17239.166
$ws.Range("M12").Value = -17069.166
$ws.Range("H32").Value = 2854.5833
$ws.Range("I32").Value = 2737.8
$ws.Range("J32").Value = 2938
$ws.Range("K32").Value = 2737.8
$ws.Range("L32").Value = 2938
$ws.Range("M32").Value = -2411.8
$ws.Range("N32").Value = -3590
$ws.Range("H106").Value = 1097.8
$ws.Range("I106").Value = 1122.25
$ws.Range("K106").Value = 1122.25
$ws.Range("M106").Value = -491.25
$ws.Range("H137").Value = 1378.7858
$ws.Range("I137").Value = 1445.3043
$ws.Range("J137").Value = 1072.8
$ws.Range("K137").Value = 4335.9129
$ws.Range("L137").Value = 3218.4
$ws.Range("M137").Value = -1785.9129
$ws.Range("N137").Value = -8318.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5206
$ws.Range("I122").Value = 2378.3572
$ws.Range("K122").Value = 7135.071599999999
$ws.Range("M122").Value = -4685.071599999999
$ws.Range("H132").Value = 5884054
$ws.Range("I132").Value = 6668442.5
$ws.Range("J132").Value = 1144
$ws.Range("K132").Value = 20005327.5
$ws.Range("L132").Value = 3432
$ws.Range("M132").Value = -20002797.5
$ws.Range("N132").Value = -8492

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3175.5
$ws.Range("I86").Value = 3523
$ws.Range("J86").Value = 2133
$ws.Range("K86").Value = 3523
$ws.Range("L86").Value = 2133
$ws.Range("M86").Value = -2400
$ws.Range("N86").Value = -4379
$ws.Range("H89").Value = 3175.5
$ws.Range("I89").Value = 3523
$ws.Range("J89").Value = 2133
$ws.Range("K89").Value = 17615
$ws.Range("L89").Value = 10665
$ws.Range("M89").Value = -11999
$ws.Range("N89").Value = -21897
$ws.Range("H99").Value = 2801.6924
$ws.Range("I99").Value = 2724.7778
$ws.Range("K99").Value = 2724.7778
$ws.Range("M99").Value = -1226.7778
$ws.Range("H105").Value = 1896.9524
$ws.Range("I105").Value = 1616.9286
$ws.Range("K105").Value = 1616.9286
$ws.Range("M105").Value = 130.0714
$ws.Range("H134").Value = 23186382
$ws.Range("I134").Value = 28338034
$ws.Range("J134").Value = 3946.5
$ws.Range("K134").Value = 85014102
$ws.Range("L134").Value = 11839.5
$ws.Range("M134").Value = -85011567
$ws.Range("N134").Value = -16909.5
$ws.Range("H135").Value = 80833.336
$ws.Range("J135").Value = 80833.336
$ws.Range("L135").Value = 80833.336
$ws.Range("N135").Value = -90973.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1597.6
$ws.Range("I16").Value = 1497
$ws.Range("K16").Value = 1497
$ws.Range("M16").Value = -1210
$ws.Range("H31").Value = 2128.946
$ws.Range("I31").Value = 2006.091
$ws.Range("K31").Value = 2006.091
$ws.Range("M31").Value = -1711.091
$ws.Range("H34").Value = 2128.946
$ws.Range("I34").Value = 2006.091
$ws.Range("K34").Value = 2006.091
$ws.Range("M34").Value = -1804.091
$ws.Range("H58").Value = 14715048
$ws.Range("J58").Value = 7286.5
$ws.Range("L58").Value = 7286.5
$ws.Range("N58").Value = -7692.5
$ws.Range("H86").Value = 10371.777
$ws.Range("I86").Value = 8877.429
$ws.Range("K86").Value = 8877.429
$ws.Range("M86").Value = -7754.429
$ws.Range("H89").Value = 10371.777
$ws.Range("I89").Value = 8877.429
$ws.Range("K89").Value = 44387.145
$ws.Range("M89").Value = -38771.145
$ws.Range("H94").Value = 2055.2222
$ws.Range("I94").Value = 2549
$ws.Range("J94").Value = 1660.2
$ws.Range("K94").Value = 2549
$ws.Range("L94").Value = 1660.2
$ws.Range("M94").Value = -2098
$ws.Range("N94").Value = -2562.2
$ws.Range("H107").Value = 91675.37
$ws.Range("I107").Value = 304.125
$ws.Range("K107").Value = 304.125
$ws.Range("M107").Value = 1615.875
$ws.Range("H113").Value = 1597.6
$ws.Range("I113").Value = 1497
$ws.Range("K113").Value = 1497
$ws.Range("M113").Value = 673
$ws.Range("H132").Value = 28574984
$ws.Range("I132").Value = 30306608
$ws.Range("K132").Value = 90919824
$ws.Range("M132").Value = -90917294
$ws.Range("H134").Value = 19233710
$ws.Range("I134").Value = 25002592
$ws.Range("K134").Value = 75007776
$ws.Range("M134").Value = -75005241
$ws.Range("H136").Value = 14715048
$ws.Range("J136").Value = 7286.5
$ws.Range("L136").Value = 21859.5
$ws.Range("N136").Value = -26959.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 689.5
$ws.Range("J31").Value = 879
$ws.Range("L31").Value = 2637
$ws.Range("N31").Value = -3213
$ws.Range("H133").Value = 16368.75
$ws.Range("J133").Value = 17300
$ws.Range("L133").Value = 51900
$ws.Range("N133").Value = -62020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 77.5
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 76.666664
$ws.Range("K3").Value = 80
$ws.Range("L3").Value = 76.666664
$ws.Range("M3").Value = 36
$ws.Range("N3").Value = -308.666664
$ws.Range("H11").Value = 35377.625
$ws.Range("I11").Value = 35377.625
$ws.Range("K11").Value = 35377.625
$ws.Range("M11").Value = -35238.625
$ws.Range("H102").Value = 2864.5789
$ws.Range("I102").Value = 2585.8235
$ws.Range("K102").Value = 2585.8235
$ws.Range("M102").Value = -963.8235
$ws.Range("H122").Value = 4222.8696
$ws.Range("I122").Value = 2964.375
$ws.Range("J122").Value = 7099.4287
$ws.Range("K122").Value = 8893.125
$ws.Range("L122").Value = 21298.2861
$ws.Range("M122").Value = -6443.125
$ws.Range("N122").Value = -26198.2861
$ws.Range("H126").Value = 2583.0833
$ws.Range("I126").Value = 2777.5557
$ws.Range("K126").Value = 8332.667099999999
$ws.Range("M126").Value = -5862.667099999999
$ws.Range("H132").Value = 4810348
$ws.Range("I132").Value = 5002682
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 15008046
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -15005516
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4253.778
$ws.Range("I7").Value = 4144.875
$ws.Range("J7").Value = 5125
$ws.Range("K7").Value = 4144.875
$ws.Range("L7").Value = 5125
$ws.Range("M7").Value = -4032.875
$ws.Range("N7").Value = -5349
$ws.Range("H46").Value = 2770.5833
$ws.Range("I46").Value = 2830
$ws.Range("K46").Value = 2830
$ws.Range("M46").Value = -2642
$ws.Range("H96").Value = 39993.5
$ws.Range("J96").Value = 39993.5
$ws.Range("L96").Value = 39993.5
$ws.Range("N96").Value = -45485.5
$ws.Range("H100").Value = 19445774
$ws.Range("I100").Value = 19445774
$ws.Range("K100").Value = 19445774
$ws.Range("M100").Value = -19445233
$ws.Range("H122").Value = 11541
$ws.Range("I122").Value = 14749
$ws.Range("J122").Value = 5125
$ws.Range("K122").Value = 44247
$ws.Range("L122").Value = 15375
$ws.Range("M122").Value = -41797
$ws.Range("N122").Value = -20275
$ws.Range("H126").Value = 4253.778
$ws.Range("I126").Value = 4144.875
$ws.Range("J126").Value = 5125
$ws.Range("K126").Value = 12434.625
$ws.Range("L126").Value = 15375
$ws.Range("M126").Value = -9964.625
$ws.Range("N126").Value = -20315
$ws.Range("H132").Value = 10671134
$ws.Range("I132").Value = 15488409
$ws.Range("J132").Value = 4312.143
$ws.Range("K132").Value = 46465227
$ws.Range("L132").Value = 12936.429
$ws.Range("M132").Value = -46462697
$ws.Range("N132").Value = -17996.429
$ws.Range("H136").Value = 3224.1875
$ws.Range("I136").Value = 3781.4546
$ws.Range("J136").Value = 1998.2
$ws.Range("K136").Value = 11344.3638
$ws.Range("L136").Value = 5994.6
$ws.Range("M136").Value = -8794.363799999999
$ws.Range("N136").Value = -11094.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 15002.5
$ws.Range("I19").Value = 20005
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 20005
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = -19831
$ws.Range("N19").Value = -10348
$ws.Range("H41").Value = 36199.4
$ws.Range("J41").Value = 38249.25
$ws.Range("L41").Value = 38249.25
$ws.Range("N41").Value = -39029.25
$ws.Range("H81").Value = 1523.0526
$ws.Range("I81").Value = 1523.0526
$ws.Range("K81").Value = 3046.1052
$ws.Range("M81").Value = -1985.1052
$ws.Range("H84").Value = 1523.0526
$ws.Range("I84").Value = 1523.0526
$ws.Range("K84").Value = 15230.526
$ws.Range("M84").Value = -9926.526
$ws.Range("H107").Value = 845.86664
$ws.Range("J107").Value = 740
$ws.Range("L107").Value = 2220
$ws.Range("N107").Value = -6060
$ws.Range("H123").Value = 111333.336
$ws.Range("I123").Value = 89000
$ws.Range("J123").Value = 122500
$ws.Range("K123").Value = 89000
$ws.Range("L123").Value = 122500
$ws.Range("M123").Value = -84100
$ws.Range("N123").Value = -132300
$ws.Range("H132").Value = 19236562
$ws.Range("I132").Value = 27779146
$ws.Range("J132").Value = 15747.5
$ws.Range("K132").Value = 83337438
$ws.Range("L132").Value = 47242.5
$ws.Range("M132").Value = -83334908
$ws.Range("N132").Value = -52302.5
$ws.Range("H136").Value = 26088496
$ws.Range("I136").Value = 30001422
$ws.Range("K136").Value = 90004266
$ws.Range("M136").Value = -90001716
